$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C54").Value = "comoment_test1"
$ws.Range("A54").Value = "Comoment1"
$ws.Range("B54").Value = "Test coSkewness"

$ws.Range("A55").Value = "Comoment2"
$ws.Range("B55").Value = "Test coKurtosis"
$ws.Range("C55").Value = "comoment_test2"

$ws.Range("D55").Select()
